# Applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.217.66"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.704.15"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.75"
$ws.Range("E5").Value = "  -2.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5297"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2642"
$ws.Range("E8").Value = "  -4.86%  "
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.80"
$ws.Range("E10").Value = "  -4.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07637"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.574"
$ws.Range("E12").Value = "  -3.03%  "
$ws.Range("D13").Value = "1.705.61"
$ws.Range("E13").Value = "  -1.70%  "
$ws.Range("D14").Value = "1.942.26"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5724"
$ws.Range("E15").Value = "  -4.92%  "
$ws.Range("D16").Value = "0.0₅8163"
$ws.Range("E16").Value = "  -2.75%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.41"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").Value = "27.211.18"
$ws.Range("E18").Value = "  -2.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.01"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.659"
$ws.Range("E21").Value = "  -3.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.44"
$ws.Range("E22").Value = "  -4.61%  "
$ws.Range("E23").Value = "  -4.73%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.23"
$ws.Range("E25").Value = "  -2.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.747"
$ws.Range("E26").Value = "  +6.42%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.249"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.27"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05361"
$ws.Range("E30").Value = "  -4.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.495"
$ws.Range("E32").Value = "  -6.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.408"
$ws.Range("E33").Value = "  -3.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.635"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.868"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("E36").Value = "  -0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9456"
$ws.Range("E37").Value = "  -4.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5843"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01628"
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.868"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "1.038.47"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8375"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.91"
$ws.Range("E44").Value = "  -1.98%  "
$ws.Range("D45").Value = "1.847.36"
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "0.0₈114"
$ws.Range("E46").Value = "  +2.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.88"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4486"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.002"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06567"
$ws.Range("E50").Value = "  +10.66%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.061"
$ws.Range("E51").Value = "  -2.72%  "

Write-Host "Applied 96 cell updates"
